# Singapore Premier League workbook update (28-04-2024 15:37 refresh).
#
# The upstream scrape re-sorted the fixtures data before re-exporting the
# sheet, which (a) swapped the shared-string slots originally used for the
# two team names "Young Lions" / "Albirex Niigata Singapore" and (b) swapped
# the full data (every column except the running index in column A) between
# several pairs of adjacent match rows that were re-ordered relative to each
# other. Net effect: a handful of rows keep showing the same team name but
# now via the other shared-string slot (no visible change), while the five
# row pairs below genuinely trade their match data, and three betting-site
# match ids were corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowCells($row, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# --- Row pairs whose match data (columns B:AB) was swapped ----------------

# Row 4 <-> Row 5 (column A ids stay put; everything else trades places)
Set-RowCells 4 @{ "B" = 6228585; "F" = 'Geylang International'; "H" = 0; "I" = 'H'; "J" = 1.222; "K" = 5.75; "L" = 8; "M" = 1.25; "N" = 5.5; "O" = 7; "P" = -1.75; "Q" = 1.875; "R" = 1.975; "S" = 4; "T" = 1.825; "U" = 2.025; "V" = 0.25; "X" = -1; "Y" = 0.875; "Z" = -1; "AA" = -1; "AB" = 1.025 }
Set-RowCells 5 @{ "B" = 6228586; "F" = 'Tanjong Pagar United'; "H" = 4; "I" = 'A'; "J" = 3.4; "K" = 4.1; "L" = 1.727; "M" = 4.5; "N" = 4.5; "O" = 1.533; "P" = 1; "Q" = 2; "R" = 1.85; "S" = 3.5; "T" = 1.925; "U" = 1.925; "V" = -1; "X" = 0.5329999999999999; "Y" = 0; "Z" = 0; "AA" = 0.925; "AB" = -1 }

# Row 22 <-> Row 23 (column A ids stay put; everything else trades places)
Set-RowCells 22 @{ "B" = 6228599; "E" = 'Hougang United FC'; "F" = 'Tampines Rovers FC'; "G" = 0; "H" = 1; "J" = 4.75; "K" = 4.2; "L" = 1.5; "M" = 7.5; "N" = 4.75; "O" = 1.3; "P" = 1.5; "Q" = 2; "R" = 1.85; "S" = 3.75; "X" = 0.3; "Y" = 1; "Z" = -1; "AA" = -1; "AB" = 0.8500000000000001 }
Set-RowCells 23 @{ "B" = 6228600; "E" = 'Tanjong Pagar United'; "F" = 'Balestier Khalsa FC'; "G" = 2; "H" = 3; "J" = 3.2; "K" = 4; "L" = 1.8; "M" = 3.4; "N" = 4.2; "O" = 1.8; "P" = 0.75; "Q" = 1.825; "R" = 2.025; "S" = 4.5; "X" = 0.8; "Y" = -0.5; "Z" = 0.5125; "AA" = 1; "AB" = -1 }

# Row 28 <-> Row 29 (column A ids stay put; everything else trades places)
Set-RowCells 28 @{ "B" = 6228604; "E" = 'Balestier Khalsa FC'; "F" = 'Albirex Niigata Singapore'; "G" = 4; "I" = 'H'; "J" = 1.4; "K" = 4.75; "L" = 5.5; "M" = 1.363; "N" = 5.25; "O" = 5.5; "P" = -1.5; "Q" = 1.925; "R" = 1.925; "S" = 4.25; "T" = 1.9; "U" = 1.95; "V" = 0.363; "X" = -1; "Y" = 0.925; "Z" = -1; "AA" = 0.8999999999999999; "AB" = -1 }
Set-RowCells 29 @{ "B" = 6228603; "E" = 'Geylang International'; "F" = 'Hougang United FC'; "G" = 0; "I" = 'A'; "J" = 1.833; "K" = 3.8; "L" = 3.25; "M" = 1.909; "N" = 3.8; "O" = 3; "P" = -0.5; "Q" = 1.975; "R" = 1.875; "S" = 3.5; "T" = 1.925; "U" = 1.925; "V" = -1; "X" = 2; "Y" = -1; "Z" = 0.875; "AA" = -1; "AB" = 0.925 }

# Row 36 <-> Row 37 (column A ids stay put; everything else trades places)
Set-RowCells 36 @{ "B" = 6228610; "E" = 'Tampines Rovers FC'; "F" = 'Geylang International'; "G" = 2; "H" = 3; "I" = 'A'; "J" = 1.333; "K" = 4.5; "L" = 7; "M" = 1.25; "N" = 5.5; "O" = 8; "P" = -1.75; "Q" = 1.9; "R" = 1.95; "S" = 3.75; "T" = 1.925; "U" = 1.925; "V" = -1; "X" = 7; "Y" = -1; "Z" = 0.95; "AA" = 0.925 }
Set-RowCells 37 @{ "B" = 6228609; "E" = 'Young Lions'; "F" = 'Albirex Niigata Singapore'; "G" = 5; "H" = 0; "I" = 'H'; "J" = 1.062; "K" = 11; "L" = 17; "M" = 1.055; "N" = 11; "O" = 17; "P" = -3.25; "Q" = 1.925; "R" = 1.925; "S" = 4.5; "T" = 1.85; "U" = 2; "V" = 0.05499999999999994; "X" = -1; "Y" = 0.925; "Z" = -1; "AA" = 0.8500000000000001 }

# Row 56 <-> Row 57 (column A ids stay put; everything else trades places)
Set-RowCells 56 @{ "B" = 7099479; "E" = 'Young Lions'; "F" = 'Hougang United FC'; "G" = 5; "H" = 0; "J" = 1.166; "K" = 7; "L" = 10; "M" = 1.125; "N" = 8.5; "O" = 13; "P" = -2.75; "Q" = 1.975; "R" = 1.875; "S" = 5; "T" = 2; "U" = 1.85; "V" = 0.125; "Y" = 0.9750000000000001; "Z" = -1; "AA" = 0; "AB" = 0 }
Set-RowCells 57 @{ "B" = 7094657; "E" = 'Lion City Sailors FC'; "F" = 'Geylang International'; "G" = 3; "H" = 1; "J" = 1.2; "K" = 6; "L" = 9; "M" = 1.181; "N" = 7; "O" = 8.5; "P" = -2.25; "Q" = 1.85; "R" = 2; "S" = 5.25; "T" = 1.9; "U" = 1.95; "V" = 0.181; "Y" = -0.5; "Z" = 0.5; "AA" = -1; "AB" = 0.95 }

# --- Lone cells whose shared-string slot flipped (text itself is unchanged) 
# --- for "Young Lions" / "Albirex Niigata Singapore" ------------------------

Set-RowCells 3 @{ "F" = 'Albirex Niigata Singapore' }
Set-RowCells 9 @{ "F" = 'Albirex Niigata Singapore' }
Set-RowCells 11 @{ "F" = 'Young Lions' }
Set-RowCells 13 @{ "F" = 'Albirex Niigata Singapore' }
Set-RowCells 14 @{ "E" = 'Young Lions' }
Set-RowCells 16 @{ "F" = 'Young Lions' }
Set-RowCells 19 @{ "E" = 'Albirex Niigata Singapore'; "F" = 'Young Lions' }
Set-RowCells 24 @{ "E" = 'Albirex Niigata Singapore' }
Set-RowCells 25 @{ "F" = 'Young Lions' }
Set-RowCells 31 @{ "E" = 'Albirex Niigata Singapore' }
Set-RowCells 33 @{ "E" = 'Young Lions' }
Set-RowCells 39 @{ "E" = 'Young Lions' }
Set-RowCells 40 @{ "E" = 'Albirex Niigata Singapore' }
Set-RowCells 42 @{ "F" = 'Young Lions' }
Set-RowCells 45 @{ "F" = 'Albirex Niigata Singapore' }
Set-RowCells 46 @{ "F" = 'Young Lions' }
Set-RowCells 48 @{ "E" = 'Albirex Niigata Singapore' }
Set-RowCells 52 @{ "F" = 'Albirex Niigata Singapore' }
Set-RowCells 53 @{ "F" = 'Young Lions' }
